# This document is a two-digit-division answer sheet: a centered date
# heading followed by a 5-column table of "a÷b=q, r" answer cells.
# The edit bumps the date by one day and swaps in a new set of division
# problems/answers for 25 of the populated cells (every 4th table row is
# blank and is left untouched).
#
# We replace each old answer with its new value using Find/Replace,
# walking the document from the top for every pair and replacing only
# the FIRST remaining match (wdReplaceOne = 1). Because the pairs below
# are listed in document order, this naturally disambiguates the two
# cells that happen to share the old text "18÷3=6, 0": the earlier
# one in the document becomes "81÷6=13, 3" and the later one becomes
# "27÷8=3, 3", matching the diff exactly. (wdReplaceAll would instead
# stomp every occurrence with the same replacement, which is wrong here.)

$d = $word.ActiveDocument

$wdReplaceOne = 1
$wdFindContinue = 1

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("2023-12-18 Monday", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "2023-12-19 Tuesday", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 2023-12-18 Monday"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("49÷7=7, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "78÷9=8, 6", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 49÷7=7, 0"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("92÷4=23, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "60÷8=7, 4", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 92÷4=23, 0"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("84÷9=9, 3", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "53÷3=17, 2", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 84÷9=9, 3"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("12÷9=1, 3", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "30÷6=5, 0", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 12÷9=1, 3"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("10÷5=2, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "48÷3=16, 0", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 10÷5=2, 0"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("30÷5=6, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "55÷8=6, 7", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 30÷5=6, 0"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("73÷8=9, 1", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "84÷5=16, 4", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 73÷8=9, 1"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("68÷8=8, 4", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "72÷9=8, 0", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 68÷8=8, 4"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("14÷9=1, 5", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "26÷7=3, 5", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 14÷9=1, 5"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("55÷7=7, 6", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "96÷7=13, 5", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 55÷7=7, 6"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("70÷2=35, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "82÷7=11, 5", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 70÷2=35, 0"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("69÷4=17, 1", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "73÷4=18, 1", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 69÷4=17, 1"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("57÷3=19, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "88÷4=22, 0", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 57÷3=19, 0"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("36÷6=6, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "20÷8=2, 4", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 36÷6=6, 0"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("92÷7=13, 1", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "48÷4=12, 0", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 92÷7=13, 1"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("85÷6=14, 1", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "29÷7=4, 1", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 85÷6=14, 1"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("78÷8=9, 6", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "51÷8=6, 3", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 78÷8=9, 6"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("62÷2=31, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "19÷9=2, 1", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 62÷2=31, 0"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("18÷3=6, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "81÷6=13, 3", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 18÷3=6, 0"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("76÷5=15, 1", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "30÷8=3, 6", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 76÷5=15, 1"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("50÷4=12, 2", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "88÷9=9, 7", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 50÷4=12, 2"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("88÷5=17, 3", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "30÷2=15, 0", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 88÷5=17, 3"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("20÷2=10, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "52÷7=7, 3", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 20÷2=10, 0"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("18÷3=6, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "27÷8=3, 3", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 18÷3=6, 0"
}

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("69÷3=23, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "63÷5=12, 3", $wdReplaceOne)
if (-not $found) {
    Write-Output "WARNING: could not find 69÷3=23, 0"
}

